$d = $word.ActiveDocument

# Remove the block of paragraphs between the empty paragraph that follows
# the title ("TBV Tags") and the section break: these paragraphs contained
# the PUMP/ACE tag markers and their descriptive text, which were deleted
# from the report.
$first = $d.Paragraphs.Item(3)
$last = $d.Paragraphs.Item($d.Paragraphs.Count)

$r = $d.Range($first.Range.Start, $last.Range.End)
$r.Delete()
